# Insert a new data row at row 323 (pushes existing rows 323:354 down to 324:355)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("323:323").Insert()

$ws.Range("A323").Value = 5
$ws.Range("B323").Value = "Macroferia Regional de Talca"
$ws.Range("C323").Value = "Maule"
$ws.Range("D323").Value = 44769
$ws.Range("E323").Value = 7
$ws.Range("F323").Value = 100114014
$ws.Range("G323").Value = "Betarraga"
$ws.Range("H323").Value = "Sin especificar"
$ws.Range("I323").Value = "Primera"
$ws.Range("J323").Value = 4000
$ws.Range("K323").Value = 800
$ws.Range("L323").Value = 800
$ws.Range("M323").Value = 800
$ws.Range("N323").Value = "$/paquete 5 unidades"
$ws.Range("O323").Value = "Región del Maule"
$ws.Range("P323").Value = 160
$ws.Range("Q323").Value = 5
$ws.Range("R323").Value = "Hortaliza"
